# Add team record columns (Wins / Losses / Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the existing header formatting (bold, centered, bordered) used by
# the other header cells (e.g. AC1) by copying its format only.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the team record for every player row (rows 2-55).
for ($r = 2; $r -le 55; $r++) {
    $ws.Cells.Item($r, 30).Value = 47
    $ws.Cells.Item($r, 31).Value = 114
    $ws.Cells.Item($r, 32).Value = 0
}
